$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25-113 down to 26-114)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new data record
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44914
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 100114007
$ws.Range("G25").Value = "Jengibre"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 190
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12368
$ws.Range("N25").Value = '$/caja 13 kilos'
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 951
$ws.Range("Q25").Value = 13
$ws.Range("R25").Value = "Hortaliza"
